$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title at the top of the document.
$removed = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Trim().StartsWith("Meta description:")) {
        $para.Range.Delete()
        $removed = $true
        break
    }
}

# 2. At the bottom, insert a new bold paragraph ("Play Aladdin's Lamp Slot
#    Game for Free | Review") right before the final paragraph (the one that
#    used to hold the "Prompt: ..." image-generation text).
$n = $d.Paragraphs.Count
$beforeLastIdx = $n - 1
$beforeLast = $d.Paragraphs.Item($beforeLastIdx)
$beforeLast.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item($beforeLastIdx + 1)
$newPara.Range.InsertXML("<w:p $wns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Aladdin's Lamp Slot Game for Free | Review</w:t></w:r></w:p>") | Out-Null

# 3. Replace the text of the final paragraph (still italic) with the new
#    meta-description copy.
$lastIdx = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($lastIdx)
$last.Range.InsertXML("<w:p $wns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Experience the magic of Arabian Nights with Aladdin's Lamp slot game. Play for free and check out our review for all the pros and cons.</w:t></w:r></w:p>") | Out-Null

Write-Output ("done; removedMeta=" + $removed)
